$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.885046362876892
$ws.Range("B1").Value = 7.446375370025635
$ws.Range("C1").Value = 6.915257453918457
$ws.Range("D1").Value = 2.308467626571655
$ws.Range("E1").Value = 1.492987036705017
